$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$NL = [char]10

# ---------------------------------------------------------------------------
# New rows 27-33: additional LeetCode problems appended below the existing
# table (row 26 "Longest ZigZag Path in Binary Tree" stays as-is).
# Cell values are written in the same order the original author typed them
# in (new vocabulary first, then the G30 note and the row 25 comment last)
# so the shared-string table comes out in the same append order.
# ---------------------------------------------------------------------------

# Row 27 - Lowest Common Ancestor of Binary Tree
$ws.Range("A27").Value = 236
$ws.Range("B27").Value = "Lowest Common Ancestor of Binary Tree"
$ws.Range("C27").Value = "Medium"
$ws.Range("C26").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = "Binary Tree"
$ws.Range("E27").Value = "BT"
$ws.Range("F27").Value = "DFS"
$ws.Range("H27").Value = "DFS"

# Row 28 - Binary Tree Right Side View
$ws.Range("A28").Value = 199
$ws.Range("B28").Value = "Binary Tree Right Side View"
$ws.Range("C28").Value = "Medium"
$ws.Range("C26").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = "Binary Tree"
$ws.Range("E28").Value = "BT"
$ws.Range("F28").Value = "BFS"
$ws.Range("G28").Value = "Only add the right child node into the resulting array"
$ws.Range("H28").Value = "BFS"

# Row 29 - Max Level Sum of Binary Tree
$ws.Range("A29").Value = 1161
$ws.Range("B29").Value = "Max Level Sum of Binary Tree"
$ws.Range("C29").Value = "Medium"
$ws.Range("C26").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D29").Value = "Binary Tree"
$ws.Range("E29").Value = "BT"
$ws.Range("F29").Value = "BFS"

# Row 30 - Keys and Rooms (note in G30 is filled in later, see below)
$ws.Range("A30").Value = 841
$ws.Range("B30").Value = "Keys and Rooms"
$ws.Range("C30").Value = "Medium"
$ws.Range("C26").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("D30").Value = "Graph"
$ws.Range("E30").Value = "Adjacent List"
$ws.Range("F30").Value = "DFS"
$ws.Range("H30").Value = "DFS"

# Row 31 - Number of Provinces
$ws.Range("A31").Value = 547
$ws.Range("B31").Value = "Number of Provinces"
$ws.Range("C31").Value = "Medium"
$ws.Range("C26").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("D31").Value = "Graph"
$ws.Range("E31").Value = "Adjacent Matrix"
$ws.Range("F31").Value = "DFS"
$ws.Range("H31").Value = "DFS"

# Back to row 30 - fill in the comment now
$ws.Range("G30").Value = "Iterate through all rooms and find a key then add it to visited"

# Row 32 - Min Cost Climbing Stairs
$ws.Range("A32").Value = 746
$ws.Range("B32").Value = "Min Cost Climbing Stairs"
$ws.Range("C32").Value = "Easy"
$ws.Range("C2").Copy()
$ws.Range("C32").PasteSpecial(-4122)
$ws.Range("D32").Value = "DP"
$ws.Range("E32").Value = "Hashmap"
$ws.Range("F32").Value = "DP"
$ws.Range("G32").Value = "Fibonacci Sequence with a twist"
$ws.Range("H32").Value = "DP"

# Row 33 - House Robber 1
$ws.Range("A33").Value = 198
$ws.Range("B33").Value = "House Robber 1"
$ws.Range("C33").Value = "Medium"
$ws.Range("C26").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("D33").Value = "DP"
$ws.Range("F33").Value = "DP"
$ws.Range("G33").Value = "Fibonacci Sequence"
$ws.Range("H33").Value = "DP"

# ---------------------------------------------------------------------------
# Row 25 (Path Sum 3): replace the short comment with the longer multi-line
# explanation last, wrap the text and grow the row height to fit it.
# ---------------------------------------------------------------------------
$ws.Range("G25").Value = "Handle 2 cases separately" + $NL + "Case 1: Find all the paths from the root node" + $NL + "Case 2: Find all paths from subtree" + $NL + "Add all the paths so it should be root + left  + right"
$ws.Range("G25").WrapText = $true
$ws.Rows.Item(25).RowHeight = 60

# ---------------------------------------------------------------------------
# View state: scroll so the new rows are visible and select G30.
# ---------------------------------------------------------------------------
$ws.Range("A19").Select() | Out-Null
$ws.Range("G30").Select() | Out-Null
